$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize the book window (cosmetic, matches the saved workbookView) -
try {
    $excel.ActiveWindow.Width = 16140
    $excel.ActiveWindow.Height = 4395
} catch {
}

# --- Add the three new data rows ---------------------------------------
$ws.Range("A2").Value = "Gaurav"
$ws.Range("B2").Value = "gabc@gmail.com"
$ws.Range("C2").Value = "Gauravnaukri@11"
$ws.Range("D2").Value = 9911227788

$ws.Range("A3").Value = "Neha"
$ws.Range("B3").Value = "nabc@gmail.com"
$ws.Range("C3").Value = "Nehanaukri@11"
$ws.Range("D3").Value = 9966332255

$ws.Range("A4").Value = "Saurabh"
$ws.Range("B4").Value = "sabc@gmail.com"
$ws.Range("C4").Value = "Saurabhnaukri@11"
$ws.Range("D4").Value = 9874563211

# --- Turn the email / "password" cells into mailto hyperlinks ----------
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:gabc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Gauravnaukri@11")

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:nabc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Nehanaukri@11")

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:sabc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Saurabhnaukri@11")

# --- Move the active selection (matches the final cursor position) -----
$ws.Range("B10").Select() | Out-Null
